$d = $word.ActiveDocument

# --- Locate the _GoBack bookmark's position ---
# In the original document the bookmark sits immediately after the words
# "designed upon" (right before " a multi-layered architecture..."). Find
# that anchor text dynamically instead of relying on a hard-coded offset.
$anchor = $d.Content
$found = $anchor.Find.Execute("designed upon", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPos = $anchor.End

# Clear everything before the bookmark (this keeps the bookmark itself
# anchored - it simply collapses to position 0 once its preceding text is
# gone).
$prefixRange = $d.Range(0, $bookmarkPos)
$prefixRange.Text = ""

# Clear everything after the bookmark, up to (but not including) the
# paragraph mark that ends the paragraph.
$endPos = $d.Content.End
$suffixRange = $d.Range(0, $endPos - 1)
$suffixRange.Text = ""

# The paragraph now contains nothing but the zero-width _GoBack bookmark,
# sitting at position 0.

# Insert the text that belongs AFTER the bookmark first: calling
# InsertAfter on the collapsed (0,0) range lands the new text immediately
# after the bookmark when the bookmark is the only thing at position 0.
$afterRange = $d.Range(0, 0)
$afterRange.InsertAfter(" every project iteration is capable of carrying out the system’s main functions.")

# Now insert the text that belongs BEFORE the bookmark: calling
# InsertBefore on the collapsed (0,0) range lands the new text immediately
# before the bookmark.
$beforeRange = $d.Range(0, 0)
$beforeRange.InsertBefore("The SuperRent system is a web-based application which is designed upon a multi-layered architecture. The primary functions, which are to buy and sell new or used cars, are handled by a system divided into structure of five layers where the tasks to carry out these functions. The User Interface where all users will commonly interact with the application. An Application layer which handles what primary functions needs to be carried out dependent upon what the user needs. The Business layer tracks financial transactions and figures when car sales or rentals are made. The Data access layer to handle the information flow between the application and the database. Finally there is the Database layer, which is responsible for storing information regarding customers, employees, car inventory, car pricing, and rental rates.   Utilising agile development and following the established architecture, the developers can ensure")

Write-Host $d.Content.Text
